$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3932.353
$ws.Range("I19").Value = 6350
$ws.Range("J19").Value = 1783.3334
$ws.Range("K19").Value = 6350
$ws.Range("L19").Value = 1783.3334
$ws.Range("M19").Value = -6175
$ws.Range("N19").Value = -2133.3334
$ws.Range("H28").Value = 335
$ws.Range("I28").Value = 335
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 335
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 150
$ws.Range("N28").ClearContents()
$ws.Range("H58").Value = 788.1429
$ws.Range("I58").Value = 270.06668
$ws.Range("J58").Value = 2083.3333
$ws.Range("K58").Value = 810.2000400000001
$ws.Range("L58").Value = 6249.999899999999
$ws.Range("M58").Value = -660.2000400000001
$ws.Range("N58").Value = -6549.999899999999
$ws.Range("H107").Value = 574.9474
$ws.Range("I107").Value = 410
$ws.Range("J107").Value = 694.9091
$ws.Range("K107").Value = 410
$ws.Range("L107").Value = 694.9091
$ws.Range("M107").Value = 1510
$ws.Range("N107").Value = -4534.9091
$ws.Range("H113").Value = 1887.5
$ws.Range("I113").Value = 1866.6666
$ws.Range("J113").Value = 1950
$ws.Range("K113").Value = 1866.6666
$ws.Range("L113").Value = 1950
$ws.Range("M113").Value = 1387.3334
$ws.Range("N113").Value = -8458
$ws.Range("H115").Value = 1000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 1000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 3000
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -6134
$ws.Range("H116").Value = 2026.8182
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 2026.8182
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 2026.8182
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -8910.8182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1114.2
$ws.Range("I2").Value = 939.53845
$ws.Range("J2").Value = 2249.5
$ws.Range("K2").Value = 939.53845
$ws.Range("L2").Value = 2249.5
$ws.Range("M2").Value = -826.53845
$ws.Range("N2").Value = -2475.5
$ws.Range("H32").Value = 19337.934
$ws.Range("I32").Value = 21140.074
$ws.Range("J32").Value = 5435.7144
$ws.Range("K32").Value = 21140.074
$ws.Range("L32").Value = 5435.7144
$ws.Range("M32").Value = -20853.074
$ws.Range("N32").Value = -6009.7144
$ws.Range("H45").Value = 735.3333
$ws.Range("I45").Value = 782.4
$ws.Range("J45").Value = 500
$ws.Range("K45").Value = 782.4
$ws.Range("L45").Value = 500
$ws.Range("M45").Value = -405.4
$ws.Range("N45").Value = -1254
$ws.Range("H110").Value = 1607.4736
$ws.Range("I110").Value = 1427.2667
$ws.Range("J110").Value = 2283.25
$ws.Range("K110").Value = 1427.2667
$ws.Range("L110").Value = 2283.25
$ws.Range("M110").Value = 617.7333000000001
$ws.Range("N110").Value = -6373.25
$ws.Range("H116").Value = 1114.2
$ws.Range("I116").Value = 939.53845
$ws.Range("J116").Value = 2249.5
$ws.Range("K116").Value = 939.53845
$ws.Range("L116").Value = 2249.5
$ws.Range("M116").Value = 1354.46155
$ws.Range("N116").Value = -6837.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1114.2
$ws.Range("I3").Value = 939.53845
$ws.Range("J3").Value = 2249.5
$ws.Range("K3").Value = 939.53845
$ws.Range("L3").Value = 2249.5
$ws.Range("M3").Value = -825.53845
$ws.Range("N3").Value = -2477.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 614
$ws.Range("I107").Value = 591
$ws.Range("J107").Value = 913
$ws.Range("K107").Value = 591
$ws.Range("L107").Value = 913
$ws.Range("M107").Value = 1329
$ws.Range("N107").Value = -4753

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 500
$ws.Range("J29").Value = 500
$ws.Range("L29").Value = 1500
$ws.Range("N29").Value = -2054
$ws.Range("H131").Value = 37382.766
$ws.Range("J131").Value = 46298.773
$ws.Range("L131").Value = 138896.319
$ws.Range("N131").Value = -148976.319

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 83346670
$ws.Range("I113").Value = 83346670
$ws.Range("K113").Value = 83346670
$ws.Range("M113").Value = -83344500
$ws.Range("H122").Value = 5152.737
$ws.Range("I122").Value = 10449.833
$ws.Range("J122").Value = 2707.923
$ws.Range("K122").Value = 31349.499
$ws.Range("L122").Value = 8123.768999999999
$ws.Range("M122").Value = -28899.499
$ws.Range("N122").Value = -13023.769
$ws.Range("H132").Value = 35876.266
$ws.Range("I132").Value = 45251.13
$ws.Range("J132").Value = 5073.143
$ws.Range("K132").Value = 135753.39
$ws.Range("L132").Value = 15219.429
$ws.Range("M132").Value = -133223.39
$ws.Range("N132").Value = -20279.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2635.1765
$ws.Range("I61").Value = 1733.1111
$ws.Range("J61").Value = 3650
$ws.Range("K61").Value = 1733.1111
$ws.Range("L61").Value = 3650
$ws.Range("M61").Value = -1531.1111
$ws.Range("N61").Value = -4054
$ws.Range("H113").Value = 2635.1765
$ws.Range("I113").Value = 1733.1111
$ws.Range("J113").Value = 3650
$ws.Range("K113").Value = 1733.1111
$ws.Range("L113").Value = 3650
$ws.Range("M113").Value = 436.8888999999999
$ws.Range("N113").Value = -7990
$ws.Range("H132").Value = 8891.207
$ws.Range("I132").Value = 14613.066
$ws.Range("J132").Value = 2760.6428
$ws.Range("K132").Value = 43839.198
$ws.Range("L132").Value = 8281.9284
$ws.Range("M132").Value = -41309.198
$ws.Range("N132").Value = -13341.9284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1706.091
$ws.Range("I107").Value = 1793.7
$ws.Range("J107").Value = 830
$ws.Range("K107").Value = 5381.1
$ws.Range("L107").Value = 2490
$ws.Range("M107").Value = -3461.1
$ws.Range("N107").Value = -6330
$ws.Range("H113").Value = 541.35297
$ws.Range("I113").Value = 453.84616
$ws.Range("K113").Value = 1361.53848
$ws.Range("M113").Value = 808.4615200000001
$ws.Range("H122").Value = 2606.1538
$ws.Range("I122").Value = 3471.4285
$ws.Range("J122").Value = 1596.6666
$ws.Range("K122").Value = 10414.2855
$ws.Range("L122").Value = 4789.9998
$ws.Range("M122").Value = -7964.2855
$ws.Range("N122").Value = -9689.9998
$ws.Range("H136").Value = 2114.875
$ws.Range("I136").Value = 2093.1
$ws.Range("K136").Value = 6279.299999999999
$ws.Range("M136").Value = -3729.299999999999
